$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Ifnk", "Ifnar1", "ECs", 2, 0.6666666666666666, 0.301559, 0.904677, 0.2570607550566701, 0.2570607550566701, 3, 1, 12.944643, 38.833929, 0.3748900893017936, 0.3748900893017936, 3.903573598436999, 35.13216238593299, 0.09636952941918155, 0.09636952941918157),
    @("ECs", "Ifnk", "Ifnar1", "FAPs", 2, 0.6666666666666666, 0.301559, 0.904677, 0.2570607550566701, 0.2570607550566701, 3, 1, 15.033452, 45.100356, 0.4353841324781401, 0.4353841324781401, 4.533472751668, 40.801254765012, 0.111920173834524, 0.111920173834524),
    @("ECs", "Ifnk", "Ifnar1", "MuSCs", 2, 0.6666666666666666, 0.301559, 0.904677, 0.2570607550566701, 0.2570607550566701, 3, 1, 6.551073333333334, 19.65322, 0.1897257782200662, 0.1897257782200662, 1.975535123326667, 17.77981610994, 0.04877105180296455, 0.04877105180296456),
    @("MuSCs", "Ifnk", "Ifnar1", "ECs", 3, 1, 0.8715450000000001, 2.614635, 0.7429392449433299, 0.7429392449433299, 3, 1, 12.944643, 38.833929, 0.3748900893017936, 0.3748900893017936, 11.281838883435, 101.536549950915, 0.278520559882612, 0.278520559882612),
    @("MuSCs", "Ifnk", "Ifnar1", "FAPs", 3, 1, 0.8715450000000001, 2.614635, 0.7429392449433299, 0.7429392449433299, 3, 1, 15.033452, 45.100356, 0.4353841324781401, 0.4353841324781401, 13.10232992334, 117.92096931006, 0.3234639586436161, 0.3234639586436161),
    @("MuSCs", "Ifnk", "Ifnar1", "MuSCs", 3, 1, 0.8715450000000001, 2.614635, 0.7429392449433299, 0.7429392449433299, 3, 1, 6.551073333333334, 19.65322, 0.1897257782200662, 0.1897257782200662, 5.709555208300001, 51.38599687470001, 0.1409547264171016, 0.1409547264171016)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowVals[$c]
    }
}
